$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume values scraped on 2023-04-09
# Column D holds price text that looks numeric (e.g. "28.305.92"); force
# text format so Excel does not coerce it into a number/date.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.305.92'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.08'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.09'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.014'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5109'
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3924'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08293'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.224'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.858.07'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.27'
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.199'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.012'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.95'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06700'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.61'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.939'
$ws.Range('E21').Value = '  -1.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.361.23'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.07'
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.260'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.072.81'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.63'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.63'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.395'
$ws.Range('E28').Value = '  -4.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '126.52'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1049'
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.032'
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.791'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.632'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.02431'
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.157'
$ws.Range('E35').Value = '  -6.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06467'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2166'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.257'
$ws.Range('E38').Value = '  +2.07%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6418'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.181'
$ws.Range('E40').Value = '  -2.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.950'
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.10'
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5990'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.88'
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.694'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.980'
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.200'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '120.65'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06855'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.03'
$ws.Range('E51').Value = '  -2.87%  '
